$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 9034.311
$ws.Range("I15").Value = 9034.311
$ws.Range("K15").Value = 27102.933
$ws.Range("M15").Value = -26933.933
$ws.Range("H62").Value = 79168910
$ws.Range("I62").Value = 105557384
$ws.Range("K62").Value = 105557384
$ws.Range("M62").Value = -105556760
$ws.Range("H65").Value = 79168910
$ws.Range("I65").Value = 105557384
$ws.Range("K65").Value = 527786920
$ws.Range("M65").Value = -527783800
$ws.Range("H86").Value = 4079.6
$ws.Range("J86").Value = 4399
$ws.Range("L86").Value = 4399
$ws.Range("N86").Value = -6645
$ws.Range("H88").Value = 1174.625
$ws.Range("I88").Value = 331.66666
$ws.Range("J88").Value = 1680.4
$ws.Range("K88").Value = 331.66666
$ws.Range("L88").Value = 1680.4
$ws.Range("M88").Value = 74.33334000000002
$ws.Range("N88").Value = -2492.4
$ws.Range("H89").Value = 4079.6
$ws.Range("J89").Value = 4399
$ws.Range("L89").Value = 21995
$ws.Range("N89").Value = -33227
$ws.Range("H91").Value = 1174.625
$ws.Range("I91").Value = 331.66666
$ws.Range("J91").Value = 1680.4
$ws.Range("K91").Value = 331.66666
$ws.Range("L91").Value = 1680.4
$ws.Range("M91").Value = 1072.33334
$ws.Range("N91").Value = -4488.4
$ws.Range("H92").Value = 1298.3846
$ws.Range("I92").Value = 1147.9
$ws.Range("K92").Value = 1147.9
$ws.Range("M92").Value = 100.0999999999999
$ws.Range("H100").Value = 19549.875
$ws.Range("I100").Value = 42832.418
$ws.Range("J100").Value = 9571.643
$ws.Range("K100").Value = 42832.418
$ws.Range("L100").Value = 9571.643
$ws.Range("M100").Value = -42291.418
$ws.Range("N100").Value = -10653.643
$ws.Range("H111").Value = 6253089.5
$ws.Range("I111").Value = 1802.75
$ws.Range("J111").Value = 12504376
$ws.Range("K111").Value = 5408.25
$ws.Range("L111").Value = 37513128
$ws.Range("M111").Value = -2341.25
$ws.Range("N111").Value = -37519262
$ws.Range("H113").Value = 3798.889
$ws.Range("I113").Value = 2998.75
$ws.Range("J113").Value = 4439
$ws.Range("K113").Value = 2998.75
$ws.Range("L113").Value = 4439
$ws.Range("M113").Value = 255.25
$ws.Range("N113").Value = -10947
$ws.Range("H129").Value = 1164
$ws.Range("I129").Value = 880.6429000000001
$ws.Range("K129").Value = 2641.9287
$ws.Range("M129").Value = 2358.0713
$ws.Range("I135").Value = 3704697.2
$ws.Range("J135").Value = 17729.5
$ws.Range("K135").Value = 33342274.8
$ws.Range("L135").Value = 159565.5
$ws.Range("M135").Value = -33339739.8
$ws.Range("N135").Value = -164635.5
$ws.Range("H137").Value = 8076970.5
$ws.Range("I137").Value = 19232650
$ws.Range("J137").Value = 20089.945
$ws.Range("K137").Value = 57697950
$ws.Range("L137").Value = 60269.835
$ws.Range("M137").Value = -57695400
$ws.Range("N137").Value = -65369.835
$ws.Range("H138").Value = 3254.3667
$ws.Range("I138").Value = 1556
$ws.Range("J138").Value = 4740.4375
$ws.Range("K138").Value = 4668
$ws.Range("L138").Value = 14221.3125
$ws.Range("M138").Value = 472
$ws.Range("N138").Value = -24501.3125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6063018
$ws.Range("I2").Value = 6063018
$ws.Range("K2").Value = 6063018
$ws.Range("M2").Value = -6062905
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H32").Value = 28471.244
$ws.Range("I32").Value = 28471.244
$ws.Range("K32").Value = 28471.244
$ws.Range("M32").Value = -28184.244
$ws.Range("H45").Value = 1751.125
$ws.Range("I45").Value = 1501.2
$ws.Range("K45").Value = 1501.2
$ws.Range("M45").Value = -1124.2
$ws.Range("H61").Value = 4191.913
$ws.Range("I61").Value = 3345.7
$ws.Range("J61").Value = 9833.333000000001
$ws.Range("K61").Value = 3345.7
$ws.Range("L61").Value = 9833.333000000001
$ws.Range("M61").Value = -3133.7
$ws.Range("N61").Value = -10257.333
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H102").Value = 2372.16
$ws.Range("I102").Value = 2372.16
$ws.Range("K102").Value = 2372.16
$ws.Range("M102").Value = -750.1599999999999
$ws.Range("H116").Value = 6063018
$ws.Range("I116").Value = 6063018
$ws.Range("K116").Value = 6063018
$ws.Range("M116").Value = -6060724
$ws.Range("H122").Value = 575.5
$ws.Range("I122").Value = 575.5
$ws.Range("K122").Value = 1726.5
$ws.Range("M122").Value = 723.5
$ws.Range("H129").Value = 55500
$ws.Range("J129").Value = 55500
$ws.Range("L129").Value = 55500
$ws.Range("N129").Value = -65500
$ws.Range("H132").Value = 5461.3076
$ws.Range("I132").Value = 3141.3333
$ws.Range("K132").Value = 9423.999899999999
$ws.Range("M132").Value = -6893.999899999999
$ws.Range("H136").Value = 4191.913
$ws.Range("I136").Value = 3345.7
$ws.Range("J136").Value = 9833.333000000001
$ws.Range("K136").Value = 10037.1
$ws.Range("L136").Value = 29499.999
$ws.Range("M136").Value = -7487.099999999999
$ws.Range("N136").Value = -34599.999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6063018
$ws.Range("I3").Value = 6063018
$ws.Range("K3").Value = 6063018
$ws.Range("M3").Value = -6062904
$ws.Range("H99").Value = 2005.1818
$ws.Range("I99").Value = 2045.7
$ws.Range("K99").Value = 2045.7
$ws.Range("M99").Value = -547.7
$ws.Range("H132").Value = 120000
$ws.Range("J132").Value = 120000
$ws.Range("L132").Value = 120000
$ws.Range("N132").Value = -130120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 402.8421
$ws.Range("I7").Value = 383
$ws.Range("J7").Value = 477.25
$ws.Range("K7").Value = 383
$ws.Range("L7").Value = 477.25
$ws.Range("M7").Value = -270
$ws.Range("N7").Value = -703.25
$ws.Range("H86").Value = 7042.778
$ws.Range("J86").Value = 7965.6
$ws.Range("L86").Value = 7965.6
$ws.Range("N86").Value = -10211.6
$ws.Range("H89").Value = 7042.778
$ws.Range("J89").Value = 7965.6
$ws.Range("L89").Value = 39828
$ws.Range("N89").Value = -51060
$ws.Range("H105").Value = 1174.4445
$ws.Range("I105").Value = 1233.125
$ws.Range("K105").Value = 1233.125
$ws.Range("M105").Value = 513.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 7264.3125
$ws.Range("I11").Value = 9365.666999999999
$ws.Range("K11").Value = 28097.001
$ws.Range("M11").Value = -27957.001
$ws.Range("H39").Value = 2262.1428
$ws.Range("J39").Value = 1999
$ws.Range("L39").Value = 5997
$ws.Range("N39").Value = -6585
$ws.Range("H55").Value = 1029.1538
$ws.Range("I55").Value = 228.6
$ws.Range("J55").Value = 1529.5
$ws.Range("K55").Value = 685.8
$ws.Range("L55").Value = 4588.5
$ws.Range("M55").Value = -508.8
$ws.Range("N55").Value = -4942.5
$ws.Range("H62").Value = 1681.5454
$ws.Range("I62").Value = 1473.3684
$ws.Range("K62").Value = 4420.1052
$ws.Range("M62").Value = -3734.1052
$ws.Range("H63").Value = 4745
$ws.Range("J63").Value = 4707.5
$ws.Range("L63").Value = 14122.5
$ws.Range("N63").Value = -15620.5
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 1681.5454
$ws.Range("I65").Value = 1473.3684
$ws.Range("K65").Value = 13260.3156
$ws.Range("M65").Value = -9828.3156
$ws.Range("H66").Value = 4745
$ws.Range("J66").Value = 4707.5
$ws.Range("L66").Value = 42367.5
$ws.Range("N66").Value = -49855.5
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H107").Value = 2475.5881
$ws.Range("I107").Value = 655.7143
$ws.Range("J107").Value = 3749.5
$ws.Range("K107").Value = 1967.1429
$ws.Range("L107").Value = 11248.5
$ws.Range("M107").Value = -47.14289999999983
$ws.Range("N107").Value = -15088.5
$ws.Range("H132").Value = 1944.9445
$ws.Range("I132").Value = 1370.8889
$ws.Range("J132").Value = 2519
$ws.Range("K132").Value = 12338.0001
$ws.Range("L132").Value = 22671
$ws.Range("M132").Value = -9808.000099999999
$ws.Range("N132").Value = -27731
$ws.Range("H137").Value = 3121.6667
$ws.Range("J137").Value = 7000
$ws.Range("L137").Value = 21000
$ws.Range("N137").Value = -31200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2276.5715
$ws.Range("I2").Value = 2276.5715
$ws.Range("K2").Value = 2276.5715
$ws.Range("M2").Value = -2163.5715
$ws.Range("H80").Value = 2605.9375
$ws.Range("J80").Value = 3054.6
$ws.Range("L80").Value = 3054.6
$ws.Range("N80").Value = -5050.6
$ws.Range("H83").Value = 2605.9375
$ws.Range("J83").Value = 3054.6
$ws.Range("L83").Value = 15273
$ws.Range("N83").Value = -25257
$ws.Range("H97").Value = 839.75
$ws.Range("I97").Value = 959.63635
$ws.Range("J97").Value = 576
$ws.Range("K97").Value = 959.63635
$ws.Range("L97").Value = 576
$ws.Range("M97").Value = -463.63635
$ws.Range("N97").Value = -1568
$ws.Range("H102").Value = 914.5238000000001
$ws.Range("I102").Value = 810.25
$ws.Range("K102").Value = 810.25
$ws.Range("M102").Value = 811.75
$ws.Range("H105").Value = 45499.5
$ws.Range("J105").Value = 45499.5
$ws.Range("L105").Value = 45499.5
$ws.Range("N105").Value = -52487.5
$ws.Range("H122").Value = 6636.65
$ws.Range("I122").Value = 6636.65
$ws.Range("K122").Value = 19909.95
$ws.Range("M122").Value = -17459.95
$ws.Range("H132").Value = 4126.2
$ws.Range("I132").Value = 2341.7368
$ws.Range("K132").Value = 7025.2104
$ws.Range("M132").Value = -4495.2104

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3441.697
$ws.Range("I22").Value = 2106.1667
$ws.Range("J22").Value = 5044.3335
$ws.Range("K22").Value = 2106.1667
$ws.Range("L22").Value = 5044.3335
$ws.Range("M22").Value = -1811.1667
$ws.Range("N22").Value = -5634.3335
$ws.Range("H27").Value = 3441.697
$ws.Range("I27").Value = 2106.1667
$ws.Range("J27").Value = 5044.3335
$ws.Range("K27").Value = 2106.1667
$ws.Range("L27").Value = 5044.3335
$ws.Range("M27").Value = -1999.1667
$ws.Range("N27").Value = -5258.3335
$ws.Range("H61").Value = 783.3333
$ws.Range("I61").Value = 783.3333
$ws.Range("K61").Value = 783.3333
$ws.Range("M61").Value = -581.3333
$ws.Range("H69").Value = 56032.6
$ws.Range("I69").Value = 90000
$ws.Range("J69").Value = 47540.75
$ws.Range("K69").Value = 90000
$ws.Range("L69").Value = 47540.75
$ws.Range("N69").Value = -49162.75
$ws.Range("M69").Value = -89189
$ws.Range("H72").Value = 56032.6
$ws.Range("I72").Value = 90000
$ws.Range("J72").Value = 47540.75
$ws.Range("K72").Value = 270000
$ws.Range("L72").Value = 142622.25
$ws.Range("N72").Value = -150734.25
$ws.Range("M72").Value = -265944
$ws.Range("H100").Value = 11368813
$ws.Range("I100").Value = 15628368
$ws.Range("K100").Value = 15628368
$ws.Range("M100").Value = -15627827
$ws.Range("H113").Value = 783.3333
$ws.Range("I113").Value = 783.3333
$ws.Range("K113").Value = 783.3333
$ws.Range("M113").Value = 1386.6667
$ws.Range("H136").Value = 5709.5
$ws.Range("J136").Value = 7136.8184
$ws.Range("L136").Value = 21410.4552
$ws.Range("N136").Value = -26510.4552

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 11115645
$ws.Range("I126").Value = 12504600
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 37513800
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -37511330
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 6605.727
$ws.Range("I132").Value = 3168.8
$ws.Range("K132").Value = 9506.400000000001
$ws.Range("M132").Value = -6976.400000000001

